$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.086.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.514.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +14.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.57%  "

# Row 7
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.525.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.72%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.92%  "

# Row 11
$ws.Range("E11").Value = "  +5.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.095.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +14.16%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.329.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.84%  "

# Row 17
$ws.Range("E17").Value = "  +1.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.506.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +14.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "503.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.95%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.98%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.28%  "

# Row 28
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("E29").Value = "  +11.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.45%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.55%  "

# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000108"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +19.70%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.68%  "

# Row 35
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.88%  "

# Row 37
$ws.Range("E37").Value = "  +7.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.329"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "46.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.65%  "

# Row 42
$ws.Range("E42").Value = "  +3.59%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.995.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.01%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "396.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0360"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "

# Row 50
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.19%  "
